# Adds 2 new Recruitment test cases to the "Recruitment" worksheet:
#   TC07_Recruitment_MarkInterviewFailedAndRejectCandidate (rows 13-14)
#   TC08_Recruitment_RejectCandidateWithoutShortlist (rows 15-16)
# Also inserts a new "CandidateInterviewStatus"/"Pass" column into the
# existing TC05 block (rows 9-10), pushing its later columns one to the
# right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recruitment")

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Phase 1: write every cell VALUE first (order chosen to reproduce the
# original authoring session as closely as possible), no formatting yet.
# ---------------------------------------------------------------------

# new interview-status column inserted ahead of the old TC05 columns
$ws.Range("F9").Value = "CandidateInterviewStatus"
$ws.Range("F10").Value = "Pass"

# rest of the shifted TC05 columns (reuse already-known strings)
$ws.Range("G9").Value = "CandidateInitialStatus"
$ws.Range("H9").Value = "CandidateStatus2"
$ws.Range("I9").Value = "CandidateCurrentStatus"
$ws.Range("G10").Value = "Interview Scheduled"
$ws.Range("H10").Value = "Interview Passed"
$ws.Range("I10").Value = "Job Offered"
$ws.Range("J10").Value = "TC05_Recruitment_MarkInterviewPassedAndOfferJobToCandidate"

# TC07 header row (row 13) - fields reused from earlier templates
$ws.Range("E13").Value = "VacancyName"
$ws.Range("F13").Value = "Keywords"
$ws.Range("G13").Value = "Notes"
$ws.Range("H13").Value = "InterviewTitle"
$ws.Range("I13").Value = "DateOfInterview"
$ws.Range("J13").Value = "TimeOfInterview"
$ws.Range("K13").Value = "CandidateInterviewStatus"
$ws.Range("L13").Value = "CandidateInitialStatus"
$ws.Range("M13").Value = "CandidateStatus2"
$ws.Range("N13").Value = "CandidateCurrentStatus"
$ws.Range("O13").Value = "CandidateStatus"
$ws.Range("B13").Value = "EmployeeFirstName"
$ws.Range("C13").Value = "EmployeeMiddleName"
$ws.Range("D13").Value = "EmployeeLastName"

# TC07 data row (row 14)
$ws.Range("B14").Value = "Renu"
$ws.Range("D14").Value = "Agarwal"
$ws.Range("C14").Value = "M"
$ws.Range("E14").Value = "QA Automation"
$ws.Range("F14").Value = "Selenium,JAVA,SQL,Eclipse"
$ws.Range("G14").Value = "Skilled QA Engineer"
$ws.Range("H14").Value = "Automation Tester"
$ws.Range("I14").Value = "2024-18-03"
$ws.Range("J14").Value = "02:00 PM"
$ws.Range("M14").Value = "Interview Failed"
$ws.Range("K14").Value = "Fail"
$ws.Range("N14").Value = "Rejected"
$ws.Range("L14").Value = "Interview Scheduled"
$ws.Range("O14").Value = "Application Initiated"
$ws.Range("A13").Value = "TC07_Recruitment_MarkInterviewFailedAndRejectCandidate"
$ws.Range("P14").Value = "TC07_Recruitment_MarkInterviewFailedAndRejectCandidate"

# TC08 data row (row 16) keyword field
$ws.Range("F16").Value = "Selenium,SQL,Eclipse"

# TC08 header row (row 15)
$ws.Range("J15").Value = "RejectionNote"
$ws.Range("E15").Value = "VacancyName"
$ws.Range("F15").Value = "Keywords"
$ws.Range("G15").Value = "Notes"
$ws.Range("H15").Value = "CandidateInitialStatus"
$ws.Range("I15").Value = "CandidateCurrentStatus"
$ws.Range("B15").Value = "EmployeeFirstName"
$ws.Range("C15").Value = "EmployeeMiddleName"
$ws.Range("D15").Value = "EmployeeLastName"

# TC08 data row (row 16) remainder
$ws.Range("J16").Value = "Qualifications and experience not enough to meet the requirements."
$ws.Range("B16").Value = "Kripa"
$ws.Range("C16").Value = "N"
$ws.Range("D16").Value = "Nair"
$ws.Range("E16").Value = "QA Automation"
$ws.Range("G16").Value = "QA Engineer"
$ws.Range("H16").Value = "Application Initiated"
$ws.Range("I16").Value = "Rejected"
$ws.Range("A15").Value = "TC08_Recruitment_RejectCandidateWithoutShortlist"
$ws.Range("K16").Value = "TC08_Recruitment_RejectCandidateWithoutShortlist"

# ---------------------------------------------------------------------
# Phase 2: apply cell formatting by copying format from donor cells that
# already carry the desired style.
# ---------------------------------------------------------------------

function Copy-Style {
    param([string]$DonorAddr, [string[]]$TargetAddrs)
    $ws.Range($DonorAddr).Copy()
    foreach ($addr in $TargetAddrs) {
        $ws.Range($addr).PasteSpecial($xlPasteFormats)
    }
}

# Row 9 / Row 10 (shifted TC05 columns)
Copy-Style "E9" @("F9")
Copy-Style "F9" @("G9", "H9", "I9")
Copy-Style "F2" @("F10")
Copy-Style "G2" @("I10")
Copy-Style "A1" @("J10")

# Row 13 header (TC07)
Copy-Style "A1" @("A13")
Copy-Style "B1" @("B13", "C13", "D13")
Copy-Style "H1" @("E13", "F13", "G13", "H13", "I13", "J13", "K13", "L13", "M13", "N13", "O13")

# Row 14 data (TC07)
Copy-Style "F2" @("F14", "H14", "K14")
Copy-Style "G8" @("I14")
Copy-Style "K2" @("J14")
Copy-Style "A1" @("P14")

# Row 15 header (TC08)
Copy-Style "A1" @("A15")
Copy-Style "B1" @("B15", "C15", "D15")
Copy-Style "H1" @("E15", "F15", "G15", "H15", "I15", "J15")

# Row 16 data (TC08)
Copy-Style "F2" @("F16")
Copy-Style "A1" @("K16")

$excel.CutCopyMode = $false

$ws.Range("H18").Select()

Write-Output "done"
